$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("weibull")
$ws1.Range("B2").Value = -3.58822680878078
$ws1.Range("C2").Value = 0.250481693413968
$ws1.Range("B3").Value = 0.327498708634429
$ws1.Range("C3").Value = 0.118770393218704

$ws2 = $wb.Worksheets.Item("lognormal")
$ws2.Range("B2").Value = 2.90225813456849
$ws2.Range("C2").Value = 0.237300586670161
$ws2.Range("B3").Value = -1.14817099226754
$ws2.Range("C3").Value = 0.105467977729247

$ws3 = $wb.Worksheets.Item("llogis")
$ws3.Range("B2").Value = -2.5105798322056
$ws3.Range("C2").Value = 0.127186380028304
$ws3.Range("B3").Value = 0.647579535658755
$ws3.Range("C3").Value = 0.0977150205592001

$ws4 = $wb.Worksheets.Item("gompertz")
$ws4.Range("B2").Value = -3.06419331756653
$ws4.Range("C2").Value = 0.185570645347155
$ws4.Range("B3").Value = 0.0121285462350868
$ws4.Range("C3").Value = 0.0163362747356503

$ws6 = $wb.Worksheets.Item("weibull cov")
$ws6.Range("A2").Value = 0.0627410787355289
$ws6.Range("B2").Value = -0.0232967607253938
$ws6.Range("A3").Value = -0.0232967607253938
$ws6.Range("B3").Value = 0.0141064063053257

$ws7 = $wb.Worksheets.Item("lognormal cov")
$ws7.Range("A2").Value = 0.0563115684340025
$ws7.Range("B2").Value = -0.021248020735707
$ws7.Range("A3").Value = -0.021248020735707
$ws7.Range("B3").Value = 0.0111234943262969

$ws8 = $wb.Worksheets.Item("llogis cov")
$ws8.Range("A2").Value = 0.0161763752647042
$ws8.Range("B2").Value = 0.00484636394239416
$ws8.Range("A3").Value = 0.00484636394239416
$ws8.Range("B3").Value = 0.00954822524288489

$ws9 = $wb.Worksheets.Item("gompertz cov")
$ws9.Range("A2").Value = 0.0344364644145594
$ws9.Range("B2").Value = -0.00195349293902929
$ws9.Range("A3").Value = -0.00195349293902929
$ws9.Range("B3").Value = 0.000266873872238646
